$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.637.33"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.586.89"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.02"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.51"
$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.492"
$ws.Range("E8").Value = "  +2.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.136"
$ws.Range("E9").Value = "  -0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.97"
$ws.Range("E10").Value = "  -1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.195.82"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000208"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.08"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.569.95"
$ws.Range("E15").Value = "  +0.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.709.26"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.115"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.46"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.22"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.04"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.72"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.620"
$ws.Range("E22").Value = "  +2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.26"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.727.56"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  -2.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.05"
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.51"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.20"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.584.28"
$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.44"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.45"
$ws.Range("E33").Value = "  -3.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("E34").Value = "  -2.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.85"
$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.72"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.63"
$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "173.33"
$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0853"
$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.891"
$ws.Range("E42").Value = "  -0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.94"
$ws.Range("E43").Value = "  +1.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  +5.04%  "

$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.19"
$ws.Range("E47").Value = "  -3.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.21"
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.67"
$ws.Range("E49").Value = "  +2.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.942"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.236"
$ws.Range("E51").Value = "  -1.33%  "
